# Apply the "feature" sheet status updates described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("feature")

# Mark the checkout/order/invoice rows as "done" in the status column (E).
$ws.Range("E8").Value = "done"
$ws.Range("E9").Value = "done"
$ws.Range("E11").Value = "done"
$ws.Range("E13").Value = "done"

# Flag the "view order" row with a note about the outstanding issue,
# highlighted with a yellow fill.
$ws.Range("E12").Value = "issue ,can not populate order data into db"
$ws.Range("E12").Interior.Color = 65535

# Leave the cursor on D10, matching the author's last saved selection.
$ws.Range("D10").Select() | Out-Null
